$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row formatting (style index used by row 1) down to the new rows 2-4
$headerRow = $ws.Range("A1:K1")
$headerRow.Copy()
$newRows = $ws.Range("A2:K4")
$newRows.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2: vbai-body-height
$ws.Cells.Item(2, 1).Value = "vbai-body-height"
$ws.Cells.Item(2, 2).Value = "VBAI Body Height"
$ws.Cells.Item(2, 3).Value = "null#vital-signs"
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = "LOINC#8302-2"
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(2, 7).Value = "dateTimeĵ, Periodĵ"
$ws.Cells.Item(2, 8).Value = "Quantityĵ"
$ws.Cells.Item(2, 9).Value = "optional"
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""

# Row 3: vbai-body-weight
$ws.Cells.Item(3, 1).Value = "vbai-body-weight"
$ws.Cells.Item(3, 2).Value = "VBAI Body Weight"
$ws.Cells.Item(3, 3).Value = "null#vital-signs"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = "LOINC#29463-7"
$ws.Cells.Item(3, 6).Value = ""
$ws.Cells.Item(3, 7).Value = "dateTimeĵ, Periodĵ"
$ws.Cells.Item(3, 8).Value = "Quantityĵ"
$ws.Cells.Item(3, 9).Value = "optional"
$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(3, 11).Value = ""

# Row 4: vbai-vital-signs
$ws.Cells.Item(4, 1).Value = "vbai-vital-signs"
$ws.Cells.Item(4, 2).Value = "VBAI Vital Signs"
$ws.Cells.Item(4, 3).Value = "null#vital-signs"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = "http://hl7.org/fhir/us/core/ValueSet/us-core-vital-signs (extensible)"
$ws.Cells.Item(4, 7).Value = "dateTimeĵ, Periodĵ"
$ws.Cells.Item(4, 8).Value = "Quantityĵ, CodeableConceptĵ, stringĵ, booleanĵ, integerĵ, Rangeĵ, Ratioĵ, SampledDataĵ, timeĵ, dateTimeĵ, Periodĵ"
$ws.Cells.Item(4, 9).Value = "optional"
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = ""

Write-Host "Applied vbai-fhir observations-summary update"
